$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is a new rule row; build it by copying row 3's formatting/values
# (same pattern used for row 3), then update the two cells that actually
# differ (Target Address / TCP Port).
$ws.Range("A3:H3").Copy($ws.Range("A4:H4"))

# New target address + full port list on the new row.
$ws.Range("C4").Value = "192.168.0.0/24"
$ws.Range("E4").Value = "22, 2866, 3306"

# The original rule (row 3) now only needs the trimmed port list.
$ws.Range("E3").Value = "2866, 3306"

# Data validation: D4 used to carry its own "직접입력, ALL" list validation
# separate from D3's "직접입력(Direct Input), ALL" one. Collapse them into a
# single D3:D4 validation (and re-add B3:B4 after it so the document order
# matches: D3:D4 rule first, then B3:B4 rule).
$ws.Range("D4").Validation.Delete()
$ws.Range("D3").Validation.Delete()
$ws.Range("B3:B4").Validation.Delete()
$ws.Range("D3:D4").Validation.Add(3, 1, 1, '"직접입력(Direct Input), ALL"')
$ws.Range("B3:B4").Validation.Add(3, 1, 1, '"Inbound, Outbound"')

# Leave the cursor where the author's session ended up.
$ws.Range("E9").Select()
